# Adds the Naive Bayes results rows on both worksheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Dep or Non-Dep Without LOSOCV" ---
$ws1 = $wb.Worksheets.Item("Dep or Non-Dep Without LOSOCV")

# Table 1 (without gender/age/work features) - Naive Bayes row 9
$ws1.Range("B9").Value = 0.67469879518072196
$ws1.Range("C9").Value = 0.59459459459459396
$ws1.Range("D9").Value = 0.64705882352941102
$ws1.Range("E9").Value = 0.61971830985915499
$ws1.Range("F9").Value = 0.67046818727490998

# Table 2 (with gender/age/work features) - Naive Bayes row 19
$ws1.Range("B19").Value = 0.67469879518072196
$ws1.Range("C19").Value = 0.59459459459459396
$ws1.Range("D19").Value = 0.64705882352941102
$ws1.Range("E19").Value = 0.61971830985915499
$ws1.Range("F19").Value = 0.67046818727490998
$ws1.Range("F19").NumberFormat = "0.000000"

# --- Sheet 2: "Dep or Non-Dep With LOSOCV" ---
$ws2 = $wb.Worksheets.Item("Dep or Non-Dep With LOSOCV")

# Table 1 (without gender/age/work features) - Naive Bayes row 9
$ws2.Range("B9").Value = 0.65558880513425899
$ws2.Range("C9").Value = 0.4
$ws2.Range("D9").Value = 0.304779614325068
$ws2.Range("D9").NumberFormat = "0.00000"
$ws2.Range("E9").Value = 0.33299029253574702
$ws2.Range("F9").Value = 0.65558880513425899

# Table 2 (with gender/age/work features) - Naive Bayes row 19
$ws2.Range("B19").Value = 0.65558880513425899
$ws2.Range("C19").Value = 0.4
$ws2.Range("D19").Value = 0.304779614325068
$ws2.Range("D19").NumberFormat = "0.00000"
$ws2.Range("E19").Value = 0.33299029253574702
$ws2.Range("F19").Value = 0.65558880513425899
